# Weekly refresh of the "Plátano" sub-dataset: a new week's worth of
# observations (2 rows) is inserted at the top of the data block
# (rows 672-780), pushing the existing rows down by 2 (674-782), and the
# two brand-new rows are written into the vacated 672-673 slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing block of 109 data rows (672:780) down by two rows,
#    to 674:782. Using Range.Copy with a destination preserves both values
#    and formatting (e.g. the date style on column D).
$srcBlock = $ws.Range("A672:T780")
$dstBlock = $ws.Range("A674")
$srcBlock.Copy($dstBlock)

# 2) Write the new week's two rows into the now-vacated 672:673 rows.
$ws.Range("D672").Value = 44522
$ws.Range("L672").Value = "Primera Maduro"
$ws.Range("M672").Value = 600
$ws.Range("N672").Value = 20000
$ws.Range("O672").Value = 20000
$ws.Range("P672").Value = 20000
$ws.Range("S672").Value = 1000

$ws.Range("D673").Value = 44522
$ws.Range("L673").Value = "Primera Pintón"
$ws.Range("M673").Value = 450
$ws.Range("N673").Value = 21000
$ws.Range("O673").Value = 21000
$ws.Range("P673").Value = 21000
$ws.Range("S673").Value = 1050
